# Reporte diario: refresh Desayuno breakdown for the new fecha (2025-05-15 / serial 45792),
# replacing the old 3-row sample (2025-05-05 / serial 45782) with the full 22-diet breakdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fecha = 45792

$diets = @(
    @{ Nombre = "Astringente";                  Cantidad = 3 },
    @{ Nombre = "Blanda";                        Cantidad = 19 },
    @{ Nombre = "Coronaria";                     Cantidad = 15 },
    @{ Nombre = "Hepatica";                      Cantidad = 1 },
    @{ Nombre = "Hipercalorica";                 Cantidad = 1 },
    @{ Nombre = "Hiperproteica";                 Cantidad = 2 },
    @{ Nombre = "Hipo Grasa";                    Cantidad = 6 },
    @{ Nombre = "Hipoglucida";                   Cantidad = 12 },
    @{ Nombre = "Hiposodica";                    Cantidad = 28 },
    @{ Nombre = "Liquida Clara";                 Cantidad = 4 },
    @{ Nombre = "Liquida Total";                 Cantidad = 7 },
    @{ Nombre = "Liquida Total 140 Cc";          Cantidad = 1 },
    @{ Nombre = "Liquida Total Miel 140 Cc";     Cantidad = 2 },
    @{ Nombre = "Liquida Total Nectar";          Cantidad = 9 },
    @{ Nombre = "Liquida Total Nectar 140 Cc";   Cantidad = 2 },
    @{ Nombre = "Liquida total Miel";            Cantidad = 3 },
    @{ Nombre = "Normal";                        Cantidad = 57 },
    @{ Nombre = "Renal Dialisis";                Cantidad = 4 },
    @{ Nombre = "Renal PRE Dialisis";            Cantidad = 5 },
    @{ Nombre = "Semiblanda";                    Cantidad = 28 },
    @{ Nombre = "Semiblanda Pequena";            Cantidad = 4 },
    @{ Nombre = "Todo Pure";                     Cantidad = 2 }
)

$row = 2
foreach ($diet in $diets) {
    $ws.Cells.Item($row, 1).Value = $fecha
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD"
    $ws.Cells.Item($row, 2).Value = "Desayuno"
    $ws.Cells.Item($row, 3).Value = $diet.Nombre
    $ws.Cells.Item($row, 4).Value = $diet.Cantidad
    $ws.Cells.Item($row, 5).Value = 0
    $row = $row + 1
}
